$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("metadata")

# A new "discord" contact row is being inserted right above the existing
# "portrait" row. This pushes the "portrait", "description/bio" and
# "extra message" rows (8-10) down to rows 9-11; row 14 and everything
# below is untouched, so this is a targeted shuffle of A8:C10 into
# A9:C11 rather than a full-sheet row insert. Conveniently, the per-column
# cell styles already already line up (A8/A9/A10 all share the "label"
# style, B8's hyperlink style is what B9 needs, B9/B10's plain style is
# what B10/B11 need) so only the values need to move - just the new
# discord cell needs a style no existing cell has yet.

# --- capture current contents (read everything before we overwrite it) ---
$portraitA = $ws.Cells.Item(8, 1).Value2    # "portrait"
$portraitB = $ws.Cells.Item(8, 2).Value2    # image url
$portraitUrl = $ws.Hyperlinks.Item(2).Address

$bioA = $ws.Cells.Item(9, 1).Value2         # "description/bio"
$bioB = $ws.Cells.Item(9, 2).Value2         # bio text

$extraA = $ws.Cells.Item(10, 1).Value2      # "extra message"

# grab B8's cell format (the hyperlink look) so we can stamp it onto B9
# later, since Hyperlinks.Add() applies its own (slightly different) style
$ws.Range("B8").Copy()
$ws.Range("B9").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# --- move the two hyperlinks so B9 (not B8) ends up linked ---
$opggUrl = $ws.Hyperlinks.Item(1).Address
$ws.Hyperlinks.Delete()

# --- row 11 (new row): "extra message" label moves down from row 10 ---
$ws.Cells.Item(11, 1).Value2 = $extraA
$ws.Cells.Item(11, 2).Value2 = ""
$ws.Rows.Item(11).RowHeight = 15

# --- row 10: "description/bio" moves down from row 9 (C10 text is untouched) ---
$ws.Cells.Item(10, 1).Value2 = $bioA
$ws.Cells.Item(10, 2).Value2 = $bioB

# --- row 9: "portrait" moves down from row 8 (C9 stays blank) ---
$ws.Cells.Item(9, 1).Value2 = $portraitA
$ws.Cells.Item(9, 2).Value2 = $portraitB

# --- row 8: brand new "discord" row; B8 is removed entirely, C8 untouched ---
$ws.Cells.Item(8, 1).Value2 = "discord"
$ws.Cells.Item(8, 1).Font.Bold = $true
$ws.Cells.Item(8, 1).Font.Name = "Arial"
$ws.Cells.Item(8, 1).HorizontalAlignment = -4131   # xlGeneral (no visible change, matches target xf)
$ws.Cells.Item(8, 2).Clear()

# --- recreate the hyperlinks at their (possibly new) locations ---
$ws.Hyperlinks.Add($ws.Range("B6"), $opggUrl) | Out-Null
$ws.Hyperlinks.Add($ws.Range("B9"), $portraitUrl) | Out-Null

# Hyperlinks.Add() re-stamps its own style on the cell; restore the style
# we captured from the original B8 (xlPasteFormats only, values untouched)
$ws.Range("B8").Copy()
$ws.Range("B9").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# update the selection to match the edit's final cursor position
$ws.Range("D4").Select()
